# Update "想去人数" (number of people wanting to attend) counts that changed
# between successive scrapes of the convention listing data.
#
# Sheet "展览" (exhibition list) and sheet "全部类型" (all-types combined list)
# both contain the same five events; their F-column counters each increment.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 65
$wsExhibit.Range("F5").Value = 5091
$wsExhibit.Range("F7").Value = 34
$wsExhibit.Range("F8").Value = 91
$wsExhibit.Range("F9").Value = 316

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 65
$wsAll.Range("F9").Value = 5091
$wsAll.Range("F11").Value = 34
$wsAll.Range("F12").Value = 91
$wsAll.Range("F14").Value = 316
